$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-05-08"
$ws.Range("B1").Value = "May 2022 (through May 08)"

$ws.Range("B2").Value = 3
$ws.Range("L2").Value = 3
$ws.Range("AF2").Value = 1

$ws.Range("AK4").Value = 1

$ws.Range("B8").Value = 2

$ws.Range("G12").Value = 1

$ws.Range("AA23").Value = 1

$ws.Range("B25").Value = 2

$ws.Range("G39").Value = 1

$ws.Range("G54").Value = 1
$ws.Range("V54").Value = 1

$ws.Range("B55").Value = 1

$ws.Range("L56").Value = 1

$ws.Range("G63").Value = 1

$ws.Range("AA64").Value = 1

$ws.Range("AA86").Value = 1

$ws.Range("AF92").Value = 1
